$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 8: status "offen" -> "gefixt" ------------------------------------
# Copy the formatting (named cell style "Gut"/Good, incl. date number format
# on column A) from row 7, which already carries the "gefixt" look, onto
# row 8, then update the Status cell text.
$ws.Range("A7:F7").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)
$ws.Range("E8").Value = "gefixt"

# --- Row 10: status "offen" -> "gefixt" ------------------------------------
$ws.Range("A9:F9").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)
$ws.Range("E10").Value = "gefixt"

$excel.CutCopyMode = $false

# --- Sheet view: scroll so column C is leftmost, select D7 ----------------
[void]$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
[void]$ws.Range("D7").Select()
